$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.84877743030339
$ws.Range("C2").Value = 8.467262516539574
$ws.Range("D2").Value = 6.025832147652803
$ws.Range("E2").Value = 11.08676915021225
$ws.Range("G2").Value = 57.84971770939936
$ws.Range("H2").Value = 20.84627224527631
$ws.Range("L2").Value = 10.00530866477253
$ws.Range("M2").Value = 16.65907924409281

$ws.Range("B3").Value = 18.53120132466005
$ws.Range("C3").Value = 8.103244832029704
$ws.Range("D3").Value = 5.915409022124485
$ws.Range("E3").Value = 11.08818344010213
$ws.Range("G3").Value = 57.05102847272899
$ws.Range("H3").Value = 20.76891903457282
$ws.Range("L3").Value = 10.01776881871444
$ws.Range("M3").Value = 16.61938382040234

$ws.Range("B4").Value = 18.34014054749969
$ws.Range("C4").Value = 7.869076040731149
$ws.Range("D4").Value = 5.848546960752619
$ws.Range("E4").Value = 11.08916771116861
$ws.Range("G4").Value = 56.56984056009963
$ws.Range("H4").Value = 20.72536101413633
$ws.Range("L4").Value = 10.0268943145043
$ws.Range("M4").Value = 16.59907578313822

$ws.Range("B5").Value = 18.26338709492128
$ws.Range("C5").Value = 7.771028459371651
$ws.Range("D5").Value = 5.821575146787315
$ws.Range("E5").Value = 11.08959777205918
$ws.Range("G5").Value = 56.37627279005687
$ws.Range("H5").Value = 20.70860694771635
$ws.Range("L5").Value = 10.03098385761376
$ws.Range("M5").Value = 16.59182821521824

$ws.Range("B6").Value = 18.25071237947277
$ws.Range("C6").Value = 7.75459128177686
$ws.Range("D6").Value = 5.817114246575329
$ws.Range("E6").Value = 11.0896709274059
$ws.Range("G6").Value = 56.34428906920975
$ws.Range("H6").Value = 20.7058852706308
$ws.Range("L6").Value = 10.03168531787701
$ws.Range("M6").Value = 16.59068701122449

$ws.Range("B7").Value = 18.3391007952266
$ws.Range("C7").Value = 7.867764268019882
$ws.Range("D7").Value = 5.848182044344117
$ws.Range("E7").Value = 11.08917339406588
$ws.Range("G7").Value = 56.56721956818254
$ws.Range("H7").Value = 20.725131022375
$ws.Range("L7").Value = 10.02694796614872
$ws.Range("M7").Value = 16.59897387007002

$ws.Range("B8").Value = 18.73853138287272
$ws.Range("C8").Value = 8.343994181617601
$ws.Range("D8").Value = 5.987587028729017
$ws.Range("E8").Value = 11.08723263560413
$ws.Range("G8").Value = 57.57254222472119
$ws.Range("H8").Value = 20.81878773736144
$ws.Range("L8").Value = 10.00929879200683
$ws.Range("M8").Value = 16.6445523951314

$ws.Range("B9").Value = 19.54785020075698
$ws.Range("C9").Value = 9.19116716839981
$ws.Range("D9").Value = 6.266687762943307
$ws.Range("E9").Value = 11.0843545460656
$ws.Range("G9").Value = 59.6078497589692
$ws.Range("H9").Value = 21.03336793342494
$ws.Range("L9").Value = 9.986392810712582
$ws.Range("M9").Value = 16.7658828126018

$ws.Range("B10").Value = 20.15144802176962
$ws.Range("C10").Value = 9.75843078335885
$ws.Range("D10").Value = 6.472969991205336
$ws.Range("E10").Value = 11.08281706351165
$ws.Range("G10").Value = 61.12960261159046
$ws.Range("H10").Value = 21.20937096889499
$ws.Range("L10").Value = 9.976699117263747
$ws.Range("M10").Value = 16.87403312271295

$ws.Range("B11").Value = 20.42657279396526
$ws.Range("C11").Value = 10.00416663853517
$ws.Range("D11").Value = 6.566617951206245
$ws.Range("E11").Value = 11.08224545095647
$ws.Range("G11").Value = 61.82495264392307
$ws.Range("H11").Value = 21.29329482536536
$ws.Range("L11").Value = 9.973837570713403
$ws.Range("M11").Value = 16.92724115766038

$ws.Range("B12").Value = 20.53072001384295
$ws.Range("C12").Value = 10.09542705970278
$ws.Range("D12").Value = 6.602015906197376
$ws.Range("E12").Value = 11.08204758406227
$ws.Range("G12").Value = 62.08848818379325
$ws.Range("H12").Value = 21.3256168437406
$ws.Range("L12").Value = 9.972976386890787
$ws.Range("M12").Value = 16.94795467768776

$ws.Range("B13").Value = 20.50829345483494
$ws.Range("C13").Value = 10.07585262767554
$ws.Range("D13").Value = 6.594395763667416
$ws.Range("E13").Value = 11.08208936799767
$ws.Range("G13").Value = 62.03172475201639
$ws.Range("H13").Value = 21.31863183722481
$ws.Range("L13").Value = 9.973151969867923
$ws.Range("M13").Value = 16.94346872872186

$ws.Range("B14").Value = 20.43514242357594
$ws.Range("C14").Value = 10.01171076927928
$ws.Range("D14").Value = 6.569531661080045
$ws.Range("E14").Value = 11.08222879851269
$ws.Range("G14").Value = 61.84663069711216
$ws.Range("H14").Value = 21.29594319565138
$ws.Range("L14").Value = 9.973762264324284
$ws.Range("M14").Value = 16.92893402295975

$ws.Range("B15").Value = 20.39032720829537
$ws.Range("C15").Value = 9.972187724815164
$ws.Range("D15").Value = 6.554292210145303
$ws.Range("E15").Value = 11.08231663114235
$ws.Range("G15").Value = 61.73327753089751
$ws.Range("H15").Value = 21.28211591095457
$ws.Range("L15").Value = 9.974165046086894
$ws.Range("M15").Value = 16.9201042827995

$ws.Range("B16").Value = 20.13346901393923
$ws.Range("C16").Value = 9.742121437298085
$ws.Range("D16").Value = 6.466842831460518
$ws.Range("E16").Value = 11.08285700949016
$ws.Range("G16").Value = 61.08420294998297
$ws.Range("H16").Value = 21.20396290959216
$ws.Range("L16").Value = 9.976917263229749
$ws.Range("M16").Value = 16.87063555092739

$ws.Range("B17").Value = 19.97595129618444
$ws.Range("C17").Value = 9.597810763353321
$ws.Range("D17").Value = 6.413119709332402
$ws.Range("E17").Value = 11.08322139655449
$ws.Range("G17").Value = 60.68664948665644
$ws.Range("H17").Value = 21.15699802375424
$ws.Range("L17").Value = 9.979002057324893
$ws.Range("M17").Value = 16.84130698347943

$ws.Range("B18").Value = 19.88540920102328
$ws.Range("C18").Value = 9.513648996237212
$ws.Range("D18").Value = 6.382203744325267
$ws.Range("E18").Value = 11.08344299814791
$ws.Range("G18").Value = 60.4582919380806
$ws.Range("H18").Value = 21.13034915689785
$ws.Range("L18").Value = 9.980346912352173
$ws.Range("M18").Value = 16.82481616445181

$ws.Range("B19").Value = 19.85476672617169
$ws.Range("C19").Value = 9.484955159771626
$ws.Range("D19").Value = 6.371734563930759
$ws.Range("E19").Value = 11.08352008592168
$ws.Range("G19").Value = 60.38103320994153
$ws.Range("H19").Value = 21.12138922370507
$ws.Range("L19").Value = 9.980827292663292
$ws.Range("M19").Value = 16.81929795830367

$ws.Range("B20").Value = 19.99271417726506
$ws.Range("C20").Value = 9.613292923338568
$ws.Range("D20").Value = 6.418840541050212
$ws.Range("E20").Value = 11.08318136180483
$ws.Range("G20").Value = 60.72893983151547
$ws.Range("H20").Value = 21.16195992530695
$ws.Range("L20").Value = 9.978765045516065
$ws.Range("M20").Value = 16.84439000051487

$ws.Range("B21").Value = 20.45663054068124
$ws.Range("C21").Value = 10.03059965226057
$ws.Range("D21").Value = 6.576836888036064
$ws.Range("E21").Value = 11.08218733813604
$ws.Range("G21").Value = 61.90099305626665
$ws.Range("H21").Value = 21.30259280235225
$ws.Range("L21").Value = 9.973576971514294
$ws.Range("M21").Value = 16.93318798860567

$ws.Range("B22").Value = 20.75956357844545
$ws.Range("C22").Value = 10.29286672945323
$ws.Range("D22").Value = 6.679704743751989
$ws.Range("E22").Value = 11.0816461466915
$ws.Range("G22").Value = 62.66818863172225
$ws.Range("H22").Value = 21.39765624000978
$ws.Range("L22").Value = 9.971482570370457
$ws.Range("M22").Value = 16.99450856478176

$ws.Range("B23").Value = 20.59794337659368
$ws.Range("C23").Value = 10.15385418623044
$ws.Range("D23").Value = 6.624849889112958
$ws.Range("E23").Value = 11.0819249951874
$ws.Range("G23").Value = 62.25868568159626
$ws.Range("H23").Value = 21.34663527990986
$ws.Range("L23").Value = 9.972481862993718
$ws.Range("M23").Value = 16.9614840811269

$ws.Range("B24").Value = 19.98513561719743
$ws.Range("C24").Value = 9.606297161124221
$ws.Range("D24").Value = 6.4162542431041
$ws.Range("E24").Value = 11.0831994238095
$ws.Range("G24").Value = 60.70981972278729
$ws.Range("H24").Value = 21.15971555263395
$ws.Range("L24").Value = 9.978871742913444
$ws.Range("M24").Value = 16.84299501339906

$ws.Range("B25").Value = 19.32683714246834
$ws.Range("C25").Value = 8.971531331096758
$ws.Range("D25").Value = 6.190807350445608
$ws.Range("E25").Value = 11.085032649733
$ws.Range("G25").Value = 59.05169551279194
$ws.Range("H25").Value = 20.97205506074156
$ws.Range("L25").Value = 9.991336233587125
$ws.Range("M25").Value = 16.72968370954814
